$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Rename header cells: "_old" -> "_FV2310", "_new" -> "_FV2404" ---
$oldToNew = @{
    "Segmentname_old"          = "Segmentname_FV2310"
    "Segmentgruppe_old"        = "Segmentgruppe_FV2310"
    "Segment_old"              = "Segment_FV2310"
    "Datenelement_old"         = "Datenelement_FV2310"
    "Segment ID_old"           = "Segment ID_FV2310"
    "Code_old"                 = "Code_FV2310"
    "Qualifier_old"            = "Qualifier_FV2310"
    "Beschreibung_old"         = "Beschreibung_FV2310"
    "Bedingungsausdruck_old"   = "Bedingungsausdruck_FV2310"
    "Bedingung_old"            = "Bedingung_FV2310"
    "Segmentname_new"          = "Segmentname_FV2404"
    "Segmentgruppe_new"        = "Segmentgruppe_FV2404"
    "Segment_new"              = "Segment_FV2404"
    "Datenelement_new"         = "Datenelement_FV2404"
    "Segment ID_new"           = "Segment ID_FV2404"
    "Code_new"                 = "Code_FV2404"
    "Qualifier_new"            = "Qualifier_FV2404"
    "Beschreibung_new"         = "Beschreibung_FV2404"
    "Bedingungsausdruck_new"   = "Bedingungsausdruck_FV2404"
    "Bedingung_new"            = "Bedingung_FV2404"
}

for ($col = 1; $col -le 21; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $current = $cell.Text
    if ($oldToNew.ContainsKey($current)) {
        $cell.Value = $oldToNew[$current]
    }
}

# --- 2) Turn the data range into an Excel Table (ListObject) ---
$rng = $ws.Range("A1:U83")
$tbl = $ws.ListObjects.Add(1, $rng, $null, 1, $null)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# --- 3) Freeze the header row ---
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
